$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column G is headed "K" (formerly derived from "Strike#"); regen save_data
# recalculated these values for each match row. Update per regenerated data.
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 1
$ws.Range("G7").Value = 3
$ws.Range("G8").Value = 2
$ws.Range("G10").Value = 1
